$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.505614041169197, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 3.055818435266709)
    3 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    5 = @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 9.193893886484982)
    6 = @(1.505614041169197, 1.65323645889881, 157.8057217802531, 6.48142807727062, 167.4460003575917)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    8 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    9 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
